$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 (row 8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active/selected cell in the saved view state
$ws.Range("E8").Select()
